# RU draft setup. Trading implemented.
#
# Resets the "Lv15" draft-tracking sheet: clears the rolled Qty (D) column,
# re-randomizes/re-labels the Class (C) column for the new draft, drops the
# two now-unused rows (37/38) of class data, re-applies the descending sort
# on the (now-empty) Qty column so the sortState range grows to match, and
# moves the active selection to where the user was last working.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lv15")

# New class list for this draft (row -> class name). Two rows' worth of
# entries (formerly rows 37-38) are retired entirely.
$classes = @{
    2  = "Flame Wizard"
    3  = "Reaper"
    4  = "Crossbowman"
    5  = "Aristocrat"
    6  = "MagnetMage"
    7  = "Aran"
    8  = "Lancemaster"
    9  = "Monk"
    10 = "Berserker"
    11 = "Mercedes"
    12 = "Dragon Knight"
    13 = "Saint"
    14 = "Mogall"
    15 = "Noble"
    16 = "Buccaneer"
    17 = "ConduitMage"
    18 = "Gunmaster"
    19 = "Assassin"
    20 = "General"
    21 = "ReflectorMage"
    22 = "Viking"
    23 = "Beginner"
    24 = "Crusader"
    25 = "Entombed"
    26 = "BurnMage"
    27 = "PortalMage"
    28 = "Brigand"
    29 = "Angel"
    30 = "KilnFiend"
    31 = "Golem"
    32 = "Bael"
    33 = "Kaiser"
    34 = "Jett"
    35 = "Farmer"
    36 = "Shiva"
}

# Clear the Qty column for the whole table first (the draft hasn't rolled
# quantities yet this time around) - this also zeroes the dependent % of
# Meta / Meta % summary formulas automatically on recalc.
$ws.Range("D2:D39").ClearContents()

# Write the new class names over the old ones.
foreach ($row in $classes.Keys) {
    $ws.Range("C$row").Value = $classes[$row]
}

# Rows 37 and 38 no longer hold class entries at all.
$ws.Range("C37:D38").ClearContents()

# Re-apply the descending-by-Qty sort over the extended range so the
# worksheet's stored sort state covers C2:E39 (was C2:E38).
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("D2:D39"), 0, 2) | Out-Null
$ws.Sort.SetRange($ws.Range("C2:E39"))
$ws.Sort.Header = 0
$ws.Sort.Apply()

# Move the selection/active cell to where the user left off.
$ws.Activate()
$ws.Range("I37").Select()
